# Hjemme passive tweaks - lichtwark deleted values
# Update the first four data columns (B:E) on rows 1-3 with the new
# Lichtwark-derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - sample-size / index header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 - CON
$ws.Range("B2").Value = 57.774938684272904
$ws.Range("C2").Value = 60.281544889481289
$ws.Range("D2").Value = 53.223346562156046
$ws.Range("E2").Value = 70.799075920665231

# Row 3 - STR
$ws.Range("B3").Value = 42.657777568082231
$ws.Range("C3").Value = 61.572680834779781
$ws.Range("D3").Value = 57.136207392154937
$ws.Range("E3").Value = 82.739096652589566

# Selection narrows from the full data block to just the edited columns
$excel.Application.Goto($ws.Range("B1:E3"))
